$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E4").Value = 12
$ws.Range("E6").Value = 29
$ws.Range("E12").Value = 13
$ws.Range("F14").Value = 8
$ws.Range("H14").Value = 8
$ws.Range("E15").Value = 53
$ws.Range("E16").Value = 205
$ws.Range("F16").Value = 49
$ws.Range("H16").Value = 49
